# Update automàtic: dades i banners [2026-02-05 22:49]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- DATA_EXTRACCIO timestamps (column E) ---
$ws.Range("E2").Value  = "2026-02-05 22:47:47"
$ws.Range("E3").Value  = "2026-02-05 22:47:50"
$ws.Range("E4").Value  = "2026-02-05 22:47:52"
$ws.Range("E5").Value  = "2026-02-05 22:47:55"
$ws.Range("E6").Value  = "2026-02-05 22:47:57"
$ws.Range("E7").Value  = "2026-02-05 22:47:59"
$ws.Range("E8").Value  = "2026-02-05 22:48:02"
$ws.Range("E9").Value  = "2026-02-05 22:48:04"
$ws.Range("E10").Value = "2026-02-05 22:48:06"
$ws.Range("E11").Value = "2026-02-05 22:48:09"
$ws.Range("E12").Value = "2026-02-05 22:48:11"
$ws.Range("E13").Value = "2026-02-05 22:48:13"
$ws.Range("E14").Value = "2026-02-05 22:48:16"
$ws.Range("E15").Value = "2026-02-05 22:48:18"
$ws.Range("E16").Value = "2026-02-05 22:48:21"
$ws.Range("E17").Value = "2026-02-05 22:48:23"
$ws.Range("E18").Value = "2026-02-05 22:48:26"
$ws.Range("E19").Value = "2026-02-05 22:48:29"
$ws.Range("E20").Value = "2026-02-05 22:48:31"
$ws.Range("E21").Value = "2026-02-05 22:48:34"
$ws.Range("E22").Value = "2026-02-05 22:48:36"
$ws.Range("E23").Value = "2026-02-05 22:48:39"
$ws.Range("E24").Value = "2026-02-05 22:48:41"
$ws.Range("E25").Value = "2026-02-05 22:48:44"
$ws.Range("E26").Value = "2026-02-05 22:48:46"
$ws.Range("E27").Value = "2026-02-05 22:48:49"
$ws.Range("E28").Value = "2026-02-05 22:48:51"
$ws.Range("E29").Value = "2026-02-05 22:48:54"
$ws.Range("E30").Value = "2026-02-05 22:48:56"
$ws.Range("E31").Value = "2026-02-05 22:48:59"
$ws.Range("E32").Value = "2026-02-05 22:49:01"
$ws.Range("E33").Value = "2026-02-05 22:49:04"
$ws.Range("E34").Value = "2026-02-05 22:49:06"
$ws.Range("E35").Value = "2026-02-05 22:49:09"
$ws.Range("E36").Value = "2026-02-05 22:49:11"

# --- other measurement values that carry a unit suffix (safe as plain text) ---
$ws.Range("K6").Value  = "3.7 MJ/m2"
$ws.Range("J11").Value = "994.4 hPa"
$ws.Range("M11").Value = "5.3 °C 22:08 TU"
$ws.Range("O11").Value = "1.1 °C"
$ws.Range("O12").Value = "10.7 °C"
$ws.Range("I14").Value = "8.2 mm"
$ws.Range("J15").Value = "990.2 hPa"
$ws.Range("O15").Value = "9.1 °C"
$ws.Range("I17").Value = "9.0 mm"
$ws.Range("M17").Value = "2.9 °C 22:15 TU"
$ws.Range("I18").Value = "2.5 mm"
$ws.Range("O20").Value = "-1.2 °C"
$ws.Range("K21").Value = "5.8 MJ/m2"
$ws.Range("O21").Value = "6.7 °C"
$ws.Range("K25").Value = "3.1 MJ/m2"
$ws.Range("O27").Value = "8.6 °C"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("I31").Value = "20.0 mm"
$ws.Range("M31").Value = "5.0 °C 22:29 TU"
$ws.Range("O32").Value = "12.2 °C"
$ws.Range("O34").Value = "4.7 °C"
$ws.Range("J36").Value = "992.6 hPa"

# --- bare percentage values (column H) ---
# These must stay plain text cells ("65%" etc.), but assigning a bare
# "NN%" string straight into a General-formatted cell makes the host
# reinterpret it as a numeric percentage (changing both the stored type
# and the cell style). Route the text through a throwaway helper cell
# that has been pinned to Text format first, then copy/paste-special
# just the *value* into each target — that leaves every target's
# original style untouched. The helper lives outside the sheet's used
# range (row 37) and is removed afterwards so no trace of it remains.
$helper = $ws.Range("A37")

$helper.NumberFormat = "@"
$helper.Value = "65%"
$helper.Copy()
$ws.Range("H4").PasteSpecial(-4163)

$helper.Value = "91%"
$helper.Copy()
$ws.Range("H10").PasteSpecial(-4163)

$helper.Value = "83%"
$helper.Copy()
$ws.Range("H12").PasteSpecial(-4163)

$helper.Value = "78%"
$helper.Copy()
$ws.Range("H15").PasteSpecial(-4163)

$helper.Value = "84%"
$helper.Copy()
$ws.Range("H22").PasteSpecial(-4163)

$helper.Value = "93%"
$helper.Copy()
$ws.Range("H34").PasteSpecial(-4163)

$ws.Rows(37).Delete()
